$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = "'" + '70.169.10'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  +0.42%  '
$ws.Range('D3').Value = "'" + '3.606.10'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  +2.80%  '
$ws.Range('E4').Value = '  +0.04%  '
$ws.Range('D5').Value = "'" + '603.25'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +0.49%  '
$ws.Range('D6').Value = "'" + '196.55'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +0.39%  '
$ws.Range('E7').Value = '  +0.44%  '
$ws.Range('E8').Value = '  +0.03%  '
$ws.Range('D9').Value = "'" + '0.206'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -1.39%  '
$ws.Range('D10').Value = "'" + '0.648'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -0.73%  '
$ws.Range('E11').Value = '  -0.21%  '
$ws.Range('E12').Value = '  +1.19%  '
$ws.Range('D13').Value = "'" + '9.57'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +0.34%  '
$ws.Range('D14').Value = "'" + '4.176.13'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +2.84%  '
$ws.Range('D15').Value = "'" + '13.15'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +4.43%  '
$ws.Range('D16').Value = "'" + '591.46'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -2.85%  '
$ws.Range('D17').Value = "'" + '70.290.33'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +0.32%  '
$ws.Range('B18').Value = 'WrappedEther'
$ws.Range('C18').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D18').Value = "'" + '3.630.01'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +3.30%  '
$ws.Range('B19').Value = 'Chainlink'
$ws.Range('C19').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('D19').Value = "'" + '19.17'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +0.99%  '
$ws.Range('E20').Value = '  +1.47%  '
$ws.Range('D21').Value = "'" + '0.994'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +0.15%  '
$ws.Range('D22').Value = "'" + '17.67'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -1.54%  '
$ws.Range('E23').Value = '  -0.09%  '
$ws.Range('D24').Value = "'" + '101.68'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -2.32%  '
$ws.Range('E25').Value = '  +0.19%  '
$ws.Range('E26').Value = '  -1.17%  '
$ws.Range('D27').Value = "'" + '10.77'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -1.78%  '
$ws.Range('D28').Value = "'" + '9.59'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -0.85%  '
$ws.Range('D29').Value = "'" + '33.90'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +0.67%  '
$ws.Range('E30').Value = '  +5.74%  '
$ws.Range('D31').Value = "'" + '7.13'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +0.28%  '
$ws.Range('E32').Value = '  -3.02%  '
$ws.Range('E33').Value = '  +0.87%  '
$ws.Range('D34').Value = "'" + '63.26'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +0.11%  '
$ws.Range('D35').Value = "'" + '0.0₃0892'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +7.25%  '
$ws.Range('D36').Value = "'" + '3.945.19'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +5.46%  '
$ws.Range('E37').Value = '  +1.55%  '
$ws.Range('D38').Value = "'" + '523.86'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +4.87%  '
$ws.Range('E39').Value = '  -0.01%  '
$ws.Range('D40').Value = "'" + '36.94'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +0.48%  '
$ws.Range('D41').Value = "'" + '0.391'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -0.45%  '
$ws.Range('E42').Value = '  -1.10%  '
$ws.Range('E43').Value = '  -1.55%  '
$ws.Range('E44').Value = '  -0.43%  '
$ws.Range('E45').Value = '  +3.04%  '
$ws.Range('D46').Value = "'" + '2.85'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +0.71%  '
$ws.Range('E47').Value = '  +0.43%  '
$ws.Range('E48').Value = '  -0.83%  '
$ws.Range('E49').Value = '  -0.30%  '
$ws.Range('D50').Value = "'" + '0.000254'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +3.24%  '
$ws.Range('E51').Value = '  +3.60%  '
